# markup-draft.xlsx: switch the RuleSet's Import/Variables cells over to the
# new net.cloudburo.drools.model2 package, and refresh the sheet's saved
# view state (zoom + active selection) to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -------------------------------------------------------
# B2: "Import" row value
$ws.Range("B2").Value = "net.cloudburo.drools.model2.*"

# B3: "Variables" row value
$ws.Range("B3").Value = "net.cloudburo.drools.model2.Markup markup"

# --- View-state edits ------------------------------------------------------
# Zoom out from 161% to 75%
$excel.ActiveWindow.Zoom = 75

# Move the active selection to B3 (where the variable declaration lives)
$ws.Range("B3").Select()
